# Toolroom Project Viewer - Forecasted_Hours.xlsx edit script
# 1. Rename existing sheets to "Dept." prefixed names
# 2. Create two new "Mach." sheets (Weekly / Daily) as trimmed copies of the
#    Dept sheets (12 blank task rows instead of 14 filled task rows)
# 3. Update header captions and selections to match the target state

$wb = $excel.ActiveWorkbook

$deptWeekly = $wb.Worksheets.Item(1)
$deptDaily  = $wb.Worksheets.Item(2)

$deptWeekly.Name = "Dept. Weekly Hrs"
$deptDaily.Name  = "Dept. Daily Hrs"

# --- Create the new "Mach." sheets ------------------------------------
# Copy order matters: copying Daily first then Weekly makes Excel assign
# sheetId 4 to the Daily copy and sheetId 5 to the Weekly copy (matching
# the target workbook.xml), after which we reposition the Weekly copy so
# it precedes the Daily copy in tab order.
$deptDaily.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$deptWeekly.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$machWeekly = $wb.Worksheets.Item($wb.Worksheets.Count)
$machDaily  = $wb.Worksheets.Item($wb.Worksheets.Count - 1)
$machWeekly.Move($machDaily)

# After the move, tab order is: Dept Weekly, Dept Daily, Mach Weekly, Mach Daily
$machWeekly = $wb.Worksheets.Item(3)
$machDaily  = $wb.Worksheets.Item(4)

$machWeekly.Name = "Mach. Weekly Hrs"
$machDaily.Name  = "Mach. Daily Hrs"

# Trim the copied sheets from 14 task rows down to 12 task rows by
# deleting the two rows that used to hold the last two department task
# labels (the row directly above the Total row, twice).
$machWeekly.Rows.Item(16).Delete()
$machWeekly.Rows.Item(16).Delete()
$machDaily.Rows.Item(16).Delete()
$machDaily.Rows.Item(16).Delete()

# Clear out the department task labels - the machine sheets start blank
$machWeekly.Range("B4:B15").ClearContents()
$machDaily.Range("B4:B15").ClearContents()

# --- Update header captions --------------------------------------------
$deptWeekly.Range("B1").Value = "Dept. Forecasted Hours"
$deptDaily.Range("B1").Value  = "Dept. Forecasted Hours"
$machWeekly.Range("B1").Value = "Mach. Forecasted Hours"
$machDaily.Range("B1").Value  = "Mach. Forecasted Hours"

# --- Update sheet view selections ---------------------------------------
$deptWeekly.Activate()
$deptWeekly.Range("B2:B3").Select()

$deptDaily.Activate()
$deptDaily.Range("B2:B3").Select()

$machWeekly.Activate()
$machWeekly.Range("S21").Select()

$machDaily.Activate()
$machDaily.Range("Z5").Select()
